# Logged 2021 divisional round, simulated season from conference round
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# YDS sheet: append new game-by-game yardage figures to the four long
# space-separated run strings (R/P rows for OFF/DEF columns).
# ---------------------------------------------------------------------------
$wsYDS = $wb.Worksheets.Item("YDS")

$wsYDS.Range("B2").Value = $wsYDS.Range("B2").Value2 + " 2 -1 8 3 -1 2 4 3 4 6 3 8 9 3 9 45 5 5 1 -1 9 10 0 3 0"
$wsYDS.Range("C2").Value = $wsYDS.Range("C2").Value2 + " 11 13 3 41 3 13 8 40 5 33 20 6 3 16 5"
$wsYDS.Range("B3").Value = $wsYDS.Range("B3").Value2 + " 1 1 2 5 -3 0 1 2 6 10 7 16 1 13 -1 4 2 -2"
$wsYDS.Range("C3").Value = $wsYDS.Range("C3").Value2 + " 21 8 4 57 7 19 7 7 6 15 4 22 1 5 10 11 32 7 12 11 4 7 21 8 13 7 3 19"

# ---------------------------------------------------------------------------
# OFF sheet: updated running totals for Home (row 2) and Road (row 3).
# ---------------------------------------------------------------------------
$wsOFF = $wb.Worksheets.Item("OFF")

$wsOFF.Range("C2").Value = 534
$wsOFF.Range("F2").Value = 170
$wsOFF.Range("G2").Value = 127
$wsOFF.Range("I2").Value = 8
$wsOFF.Range("J2").Value = 71
$wsOFF.Range("L2").Value = 571
$wsOFF.Range("M2").Value = 381
$wsOFF.Range("O2").Value = 42
$wsOFF.Range("Q2").Value = 1160

$wsOFF.Range("C3").Value = 278
$wsOFF.Range("E3").Value = 68
$wsOFF.Range("F3").Value = 205
$wsOFF.Range("G3").Value = 71
$wsOFF.Range("H3").Value = 55
$wsOFF.Range("I3").Value = 116
$wsOFF.Range("J3").Value = 125
$wsOFF.Range("N3").Value = 38

# ---------------------------------------------------------------------------
# DEF sheet: updated running totals for Home (row 2) and Road (row 3).
# ---------------------------------------------------------------------------
$wsDEF = $wb.Worksheets.Item("DEF")

$wsDEF.Range("B2").Value = 7
$wsDEF.Range("C2").Value = 368
$wsDEF.Range("E2").Value = 18
$wsDEF.Range("F2").Value = 115
$wsDEF.Range("G2").Value = 131
$wsDEF.Range("J2").Value = 67
$wsDEF.Range("L2").Value = 719
$wsDEF.Range("M2").Value = 476
$wsDEF.Range("Q2").Value = 1226

$wsDEF.Range("C3").Value = 429
$wsDEF.Range("E3").Value = 66
$wsDEF.Range("F3").Value = 252
$wsDEF.Range("G3").Value = 65
$wsDEF.Range("H3").Value = 53
$wsDEF.Range("I3").Value = 136
$wsDEF.Range("J3").Value = 120
$wsDEF.Range("N3").Value = 45

# ---------------------------------------------------------------------------
# ST sheet: kicking/punting counters (row 2) + appended distance strings.
# ---------------------------------------------------------------------------
$wsST = $wb.Worksheets.Item("ST")

$wsST.Range("B2").Value = 190
$wsST.Range("D2").Value = 114
$wsST.Range("F2").Value = 271
$wsST.Range("G2").Value = 258
$wsST.Range("H2").Value = 10
$wsST.Range("J2").Value = 125
$wsST.Range("K2").Value = 118

$wsST.Range("B3").Value = 109

$wsST.Range("D3").Value = $wsST.Range("D3").Value2 + " 50 27 41 48"
$wsST.Range("B4").Value = $wsST.Range("B4").Value2 + " 61 57 64"
$wsST.Range("D4").Value = $wsST.Range("D4").Value2 + " 0 0 0 0"
$wsST.Range("B5").Value = $wsST.Range("B5").Value2 + " 31 32 19"
$wsST.Range("D5").Value = $wsST.Range("D5").Value2 + " 0 5 0 0 0"
$wsST.Range("B6").Value = $wsST.Range("B6").Value2 + " 26"

# ---------------------------------------------------------------------------
# TURNS sheet: Home/Road turnover counts.
# ---------------------------------------------------------------------------
$wsTURNS = $wb.Worksheets.Item("TURNS")

$wsTURNS.Range("B2").Value = 17
$wsTURNS.Range("C2").Value = 19
$wsTURNS.Range("D3").Value = 13

# ---------------------------------------------------------------------------
# PEN sheet: offensive false-start penalty count.
# ---------------------------------------------------------------------------
$wsPEN = $wb.Worksheets.Item("PEN")

$wsPEN.Range("B2").Value = 34
